$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.436.95"
$ws.Range("E2").Value = "  -6.73%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.754.39"
$ws.Range("E3").Value = "  -5.41%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.16%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "581.22"
$ws.Range("E5").Value = "  -5.12%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "176.22"
$ws.Range("E6").Value = "  +4.32%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.748.88"
$ws.Range("E7").Value = "  -5.36%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.635"
$ws.Range("E8").Value = "  -6.61%  "

# Row 9
$ws.Range("E9").Value = "  +0.04%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.720"
$ws.Range("E10").Value = "  -7.38%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.168"
$ws.Range("E11").Value = "  -10.17%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "53.67"
$ws.Range("E12").Value = "  -4.61%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000301"
$ws.Range("E13").Value = "  -10.62%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.75"
$ws.Range("E14").Value = "  -4.53%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.354.18"
$ws.Range("E15").Value = "  -5.41%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.766.98"
$ws.Range("E16").Value = "  -5.28%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "19.54"
$ws.Range("E17").Value = "  -4.59%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.09"
$ws.Range("E18").Value = "  -7.72%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.15"
$ws.Range("E19").Value = "  -7.16%  "

# Row 20
$ws.Range("E20").Value = "  -2.70%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "68.437.10"
$ws.Range("E21").Value = "  -6.48%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "412.20"
$ws.Range("E22").Value = "  -6.36%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.56"
$ws.Range("E23").Value = "  -5.90%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "89.44"
$ws.Range("E24").Value = "  -6.22%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.10"
$ws.Range("E25").Value = "  -9.76%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.98"
$ws.Range("E26").Value = "  -8.80%  "

# Row 27
$ws.Range("B27").Value = "Toncoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.93"
$ws.Range("E27").Value = "  -6.45%  "

# Row 28
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.83"
$ws.Range("E28").Value = "  -2.41%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.98"
$ws.Range("E29").Value = "  +0.33%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.67"
$ws.Range("E30").Value = "  -8.04%  "

# Row 31
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "33.27"
$ws.Range("E31").Value = "  -7.76%  "

# Row 32
$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.88"
$ws.Range("E32").Value = "  -0.24%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "12.82"
$ws.Range("E33").Value = "  -7.75%  "

# Row 34
$ws.Range("E34").Value = "  -8.35%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "44.11"
$ws.Range("E35").Value = "  -8.54%  "

# Row 36
$ws.Range("B36").Value = "OKB"
$ws.Range("C36").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "65.61"
$ws.Range("E36").Value = "  -6.80%  "

# Row 37
$ws.Range("B37").Value = "Bittensor"
$ws.Range("C37").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "610.14"
$ws.Range("E37").Value = "  -5.94%  "

# Row 38
$ws.Range("B38").Value = "PEPE"
$ws.Range("C38").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0925"
$ws.Range("E38").Value = "  -12.51%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.406"
$ws.Range("E39").Value = "  -5.89%  "

# Row 40
$ws.Range("E40").Value = "  +0.28%  "

# Row 41
$ws.Range("E41").Value = "  -0.02%  "

# Row 42
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.137"
$ws.Range("E42").Value = "  -5.88%  "

# Row 43
$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.15"
$ws.Range("E43").Value = "  -0.61%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.08"
$ws.Range("E44").Value = "  -9.65%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0447"
$ws.Range("E45").Value = "  -7.46%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.62"
$ws.Range("E46").Value = "  +1.97%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.46"
$ws.Range("E47").Value = "  -9.50%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.136"
$ws.Range("E48").Value = "  -8.53%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.75"
$ws.Range("E49").Value = "  -13.01%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.772.71"
$ws.Range("E50").Value = "  -1.05%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.14"
$ws.Range("E51").Value = "  -8.28%  "
